$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the table row containing "Ana Cristina Soto Ruvalcaba" (row 4),
# shifting subsequent rows up. The underlying Excel Table ("FSR") and its
# autofilter/sort range are updated automatically by Excel.
$ws.Rows(4).Delete()
